$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New mouse records (rows 21-29), all dated 2017-10-06 (serial 43014)
$date = Get-Date -Year 2017 -Month 10 -Day 6 -Hour 0 -Minute 0 -Second 0

$rows = @(
    @{ Mouse = "AVI999"; Liver = "x"; Thymus = "x"; Pancreas = "x"; Spleen = "x"; MLN = "x"; Organs = "x" },
    @{ Mouse = "av567";  Liver = "x"; Thymus = "x"; Pancreas = "x"; Spleen = "x"; MLN = " "; Organs = "x" },
    @{ Mouse = "AV345";  Liver = "x"; Thymus = "x"; Pancreas = "x"; Spleen = "x"; MLN = "x"; Organs = "x" },
    @{ Mouse = "AV678";  Liver = "x"; Thymus = "x"; Pancreas = "x"; Spleen = "x"; MLN = "x"; Organs = "x" },
    @{ Mouse = "AV544";  Liver = "x"; Thymus = "x"; Pancreas = "x"; Spleen = "x"; MLN = "x"; Organs = "x" },
    @{ Mouse = "AV666";  Liver = "x"; Thymus = "x"; Pancreas = "x"; Spleen = "x"; MLN = "x"; Organs = "x" },
    @{ Mouse = "AV777";  Liver = " "; Thymus = "x"; Pancreas = "x"; Spleen = "x"; MLN = "x"; Organs = "x" },
    @{ Mouse = "AV212";  Liver = "x"; Thymus = "x"; Pancreas = "x"; Spleen = "x"; MLN = "x"; Organs = "x" },
    @{ Mouse = "AV567";  Liver = "x"; Thymus = "x"; Pancreas = "x"; Spleen = "x"; MLN = " "; Organs = "x" }
)

$r = 21
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.Mouse
    $ws.Cells.Item($r, 2).Value = $date
    $ws.Cells.Item($r, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"
    $ws.Cells.Item($r, 3).Value = $row.Liver
    $ws.Cells.Item($r, 4).Value = $row.Thymus
    $ws.Cells.Item($r, 5).Value = $row.Pancreas
    $ws.Cells.Item($r, 6).Value = $row.Spleen
    $ws.Cells.Item($r, 7).Value = $row.MLN
    $ws.Cells.Item($r, 8).Value = $row.Organs
    $r++
}
